$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C1").Value = "C"
$ws.Range("C2").Value = "-"
$ws.Range("C3").Value = 1
$ws.Range("C4").Select()
